$d = $word.ActiveDocument

# 1. Remove the old "_GoBack" bookmark that sits at the end of the first
#    paragraph (it will be re-added further down, around the pension text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2 & 3. Drop the w:hint="cs" hint from the rFonts used in the paragraph
#    marks of the "משכורת שעתית" and "אחוז משרה" paragraphs.
$d.Paragraphs(9).Range.Font.Name = $d.Paragraphs(9).Range.Font.Name
$d.Paragraphs(10).Range.Font.Name = $d.Paragraphs(10).Range.Font.Name

# 4. Strike through the "קרן השתלמות" paragraph (except for the leading
#    "-          " bullet run) and wrap it with a fresh "_GoBack" bookmark.
$p = $d.Paragraphs(21)
$full = $p.Range
$sub = $d.Range($full.Start + 11, $full.End - 1)
$sub.Font.StrikeThrough = 1

$d.Bookmarks.Add("_GoBack", $sub) | Out-Null
